$d = $word.ActiveDocument

# 1) Remove the "Meta description" paragraph that follows the title heading.
#    (Play Barbaria Free - A Stunning and Intuitive Online Slot Game / Meta description: ...)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "^Meta description:") {
        $para.Range.Delete()
        break
    }
}

# 2) Insert a new bold paragraph "Play Barbaria Free - A Stunning and Intuitive
#    Online Slot Game" right before the final (italic) paragraph of the document.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($count)
$newPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Barbaria Free - A Stunning and Intuitive Online Slot Game</w:t></w:r></w:p>')

# 3) Replace the text of the final (italic) paragraph with the meta description
#    copy, preserving its existing (italic) run formatting.
$oldText = 'Create a feature image fitting the game "Barbaria" with the following specifications: - A cartoon-style image - Features a happy Maya warrior with glasses - The environment should depict an ancient Roman setting with barbarian weapons in the background. The image should be bright and colorful, with a focus on the Maya warrior as the main character. The warrior should look happy and adventurous, holding a sword or an axe in one hand and wearing glasses. The background can be of an ancient Roman setting with barbarian weapons displayed. The image should be able to attract players looking for a thrilling game with great graphics and visuals.'
$newText = 'Read our review of the Barbaria online slot game and play for free. Enjoy stunning graphic features, two bonuses, and a medium volatility level.'
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
